$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.612369179725647
$ws.Range("B1").Value = 2.923484802246094
$ws.Range("C1").Value = 3.592628955841064
$ws.Range("D1").Value = 3.835409879684448
$ws.Range("E1").Value = 2.808867692947388
